$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 217 already exists (Serie = 04-08-2021) but is missing the BCP (B) and
# BCU (C) values - fill those in now.
$ws.Cells.Item(217, 2).Value = 187
$ws.Cells.Item(217, 3).Value = 628

# Rows 218-247 are brand new - one row per calendar day from 05-08-2021
# through 03-09-2021. Pre-format column A as Text so the "dd-mm-yyyy"
# strings are stored as shared-string text (matching the source data)
# instead of being auto-converted into Excel date serials, then restore
# the Normal style so no per-cell formatting is left behind.
$startRow = 218
$endRow = 247
$dateRange = $ws.Range("A$startRow`:A$endRow")
$dateRange.NumberFormat = "@"

$d = Get-Date -Year 2021 -Month 8 -Day 5
for ($r = $startRow; $r -le $endRow; $r++) {
    $ws.Cells.Item($r, 1).Value = $d.ToString("dd-MM-yyyy")
    $d = $d.AddDays(1)
}

$dateRange.Style = "Normal"

# Columns B/C/D/E: every new row repeats the same figures as the
# preceding rows, except the very last row (247 / 03-09-2021) which only
# carries the D and E values, same as the source diff.
for ($r = $startRow; $r -le ($endRow - 1); $r++) {
    $ws.Cells.Item($r, 2).Value = 187
    $ws.Cells.Item($r, 3).Value = 628
    $ws.Cells.Item($r, 4).Value = 3940
    $ws.Cells.Item($r, 5).Value = 30
}

$ws.Cells.Item($endRow, 4).Value = 3940
$ws.Cells.Item($endRow, 5).Value = 30

Write-Output "done"
